$wb = $excel.ActiveWorkbook

# 1. Rename the "Region" sheet to "Zone"
$wsZone = $wb.Worksheets.Item("Region")
$wsZone.Name = "Zone"

# 2. Update values on the "PV" sheet (U/V columns, rows 2-5: 4 -> 0.5, 8 -> 0)
$wsPV = $wb.Worksheets.Item("PV")
$wsPV.Range("U2").Value = 0.5
$wsPV.Range("V2").Value = 0
$wsPV.Range("U3").Value = 0.5
$wsPV.Range("V3").Value = 0
$wsPV.Range("U4").Value = 0.5
$wsPV.Range("V4").Value = 0
$wsPV.Range("U5").Value = 0.5
$wsPV.Range("V5").Value = 0

# 3. Update values on the "Slack" sheet (V2: 4 -> 2, W2: 8 -> 0)
$wsSlack = $wb.Worksheets.Item("Slack")
$wsSlack.Range("V2").Value = 2
$wsSlack.Range("W2").Value = 0

# 4. Update the selection left on the "Slack" sheet (was N8, now L3)
$wsSlack.Range("L3").Select()

# 5. Update the selection left on the "PV" sheet (was U2:V2, now L6)
$wsPV.Range("L6").Select()

# 6. Finally activate the renamed "Zone" sheet, making it the active tab
#    (keeps its own selection at J21, unchanged)
$wsZone.Activate()
